# Slide 42, "Content Placeholder 2" shape: the grammar-rule line
#   compoundStmt = "{" statement "}" .
# needs to become
#   compoundStmt = "{" statements "}" .
# with the trailing portion of the line split into three separate runs:
#   " = "  |  "{" statements "  |  "}" ."
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(42)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(2)

# Sanity check before editing (Paragraphs().Text includes a trailing CR).
if ($para.Text.TrimEnd("`r") -ne 'compoundStmt = "{" statement "}" .') {
    throw "Unexpected paragraph text: $($para.Text)"
}

# Within the paragraph (1-based), the run being split covers
# chars 13-34: ' = "{" statement "}" .'
#   chars 13-15 -> ' = '
#   chars 16-29 -> '"{" statement '  (becomes '"{" statements ')
#   chars 30-34 -> '"}" .'
# Re-assigning the text of the middle slice naturally splits the run into
# three, leaving the untouched outer slices as their own runs.
$middle = $para.Characters(16, 14)
$middle.Text = '"{" statements '

# Final text sanity check.
if ($para.Text.TrimEnd("`r") -ne 'compoundStmt = "{" statements "}" .') {
    throw "Post-edit paragraph text mismatch: $($para.Text)"
}
